$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.153.07'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.12%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.827.26'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.74%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9989'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.00'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.70%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6202'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.77%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07356'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.35%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2912'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.23%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.05'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.38%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07680'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.28%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.825.88'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.65%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.947'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.64%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6637'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.13'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.17%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008877'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -4.95%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.849'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.24%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.129.80'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.067.78'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.87%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.60'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +3.80%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.43'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.0000'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.324'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.09%  '

$ws.Range("E24").Value = '  -0.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.94'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.51%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1409'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.49%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.483'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.86%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.63'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.74%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.484'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.83%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05911'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +6.06%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.065'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.05%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.075'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.66%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.204'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.25%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.851'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7312'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.35%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.135'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.16%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.610'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.11%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.847'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.78%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.217.68'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.43%  '

$ws.Range("E40").Value = '  -2.48%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.266'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -5.02%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9151'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.42%  '

$ws.Range("E43").Value = '  +0.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.75'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.77%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.976.35'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.37%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.77'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.81%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5082'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.02%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.155'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.64%  '

$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000117'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -5.60%  '

$ws.Range("B50").Value = 'TheSandbox'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4014'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.82%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1126'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.72%  '

